$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.286.56"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.073.99"
$ws.Range("E3").Value = "  +4.90%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'234.84"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'57.23"
$ws.Range("E8").Value = "  +5.20%  "
$ws.Range("E9").Value = "  +2.87%  "
$ws.Range("D10").Value = "'58.02"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'0.0759"
$ws.Range("E11").Value = "  +0.97%  "
$ws.Range("D12").Value = "'0.102"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "2.380.50"
$ws.Range("E13").Value = "  +4.82%  "
$ws.Range("D14").Value = "'14.45"
$ws.Range("E14").Value = "  +3.84%  "
$ws.Range("D15").Value = "'20.87"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'0.773"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "'5.21"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "2.105.76"
$ws.Range("E18").Value = "  +6.30%  "
$ws.Range("D19").Value = "37.432.73"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").Value = "'6.02"
$ws.Range("E20").Value = "  +21.71%  "
$ws.Range("D21").Value = "'68.32"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "0.0₃0810"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "'223.33"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'2.42"
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'162.62"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "'8.86"
$ws.Range("E28").Value = "  +3.28%  "
$ws.Range("D29").Value = "'0.130"
$ws.Range("E29").Value = "  +7.28%  "
$ws.Range("D30").Value = "'19.19"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +6.08%  "
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").Value = "'4.44"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "'0.0620"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "'2.55"
$ws.Range("E35").Value = "  +10.38%  "
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +14.93%  "
$ws.Range("D39").Value = "'3.32"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0963"
$ws.Range("E42").Value = "  +10.06%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.468.54"
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.32"
$ws.Range("E44").Value = "  +15.91%  "
$ws.Range("D45").Value = "'95.01"
$ws.Range("E45").Value = "  +8.48%  "
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "'16.19"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("D50").Value = "'7.26"
$ws.Range("E50").Value = "  +9.66%  "
$ws.Range("D51").Value = "'2.93"
$ws.Range("E51").Value = "  +2.09%  "
